$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row data (rows 2-7) with the freshly scraped opportunity records.
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = "1328500"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328500"
$ws.Range("C2").Value = "[CC] Human Resources Talent Management Coordinator and Analyst"
$ws.Range("D2").Value = "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "176 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "DHL Group"

# Row 3
$ws.Range("A3").Value = "1328987"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328987"
$ws.Range("C3").Value = "Accelerate Romania | Architect Intern"
$ws.Range("D3").Value = "Timișoara, Romania"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "1 applicant"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Mzr Habitat"

# Row 4
$ws.Range("A4").Value = "1328980"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328980"
$ws.Range("C4").Value = "Engineering Intern"
$ws.Range("D4").Value = "Belgrade, Serbia"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "1 applicant"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "Tepma MEP Design"

# Row 5
$ws.Range("A5").Value = "1328974"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1328974"
$ws.Range("C5").Value = "Youth Volleyball Assistant Coach"
$ws.Range("D5").Value = "Belgrade, Serbia"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "OK Roda"

# Row 6
$ws.Range("A6").Value = "1328934"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1328934"
$ws.Range("C6").Value = "[CC] Global Coordinator for Talent Acquisition and Project Management (Only AIESECers. Fully read before applying)"
$ws.Range("D6").Value = "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany"
$ws.Range("E6").Value = "Yes"
$ws.Range("F6").Value = "1 applicant"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "DHL Group"

# Row 7
$ws.Range("A7").Value = "1328730"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1328730"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Bursa, Türkiye"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "4 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Orhan Holding Au"

# ---------------------------------------------------------------------------
# 2. Highlight the "PREMIUM = Yes" cells with a yellow fill (new style).
# ---------------------------------------------------------------------------
$ws.Range("E2").Interior.Color = 65535
$ws.Range("E6").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3. Remove the now-stale trailing rows (old rows 8-10 no longer scraped).
# ---------------------------------------------------------------------------
$ws.Range("A8:H10").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4. Resize columns C, D and H to fit the new (longer / shorter) content.
#    ColumnWidth setter adds Excel's standard 0.8333 padding on read-back,
#    so subtract it here to land exactly on the target stored width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 116.16666666666667
$ws.Columns.Item(4).ColumnWidth = 50.166666666666664
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668
